$wb = $excel.ActiveWorkbook

# Sheets that use the "Ano <year>" label pattern on row 1 (columns B:E)
$anoSheets = @(
    "Potencia Acumulada - SIN (MW)",
    "Geracao Periodo Medio (MWMed)",
    "Atendimento a Ponta(MW)",
    "Emissoes Totais (MtCO2eq)"
)

foreach ($sheetName in $anoSheets) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("B1").Value = "Ano " + $ws.Range("B1").Text
    $ws.Range("C1").Value = "Ano " + $ws.Range("C1").Text
    $ws.Range("D1").Value = "Ano " + $ws.Range("D1").Text
    $ws.Range("E1").Value = "Ano " + $ws.Range("E1").Text
}

# Sheet that uses the "Intervalo <range>" label pattern on row 1 (columns B:E)
$ws4 = $wb.Worksheets.Item("Potencia Incremental - SIN(MW)")
$ws4.Range("B1").Value = "Intervalo " + $ws4.Range("B1").Text
$ws4.Range("C1").Value = "Intervalo " + $ws4.Range("C1").Text
$ws4.Range("D1").Value = "Intervalo " + $ws4.Range("D1").Text
$ws4.Range("E1").Value = "Intervalo " + $ws4.Range("E1").Text

# Sheet that only has column B in row 1, also uses the "Ano <year>" pattern
$ws6 = $wb.Worksheets.Item("Custo Total (bilhões de R$)")
$ws6.Range("B1").Value = "Ano " + $ws6.Range("B1").Text
